$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value that was in C198 (becomes an empty cell)
$ws.Range("C198").Value = ""

# Append the new rows produced by the updated script run.
# The date column stores plain text (e.g. "2025-08-07"), not a real date,
# so a leading apostrophe forces Excel to keep it as text instead of
# auto-converting it to a date serial number.
$ws.Range("A199").Value = "'2025-08-07"
$ws.Range("B199").Value = "développement durable"
$ws.Range("C199").Value = 11
$ws.Range("D199").Value = 1

$ws.Range("A200").Value = "'2025-08-07"
$ws.Range("B200").Value = "développement durable"
$ws.Range("C200").Value = 13
$ws.Range("D200").Value = 1

$ws.Range("A201").Value = "'2025-08-07"
$ws.Range("B201").Value = "bonnes pratiques"
$ws.Range("C201").Value = 62
$ws.Range("D201").Value = 1

$ws.Range("A202").Value = "'2025-08-07"
$ws.Range("B202").Value = "eaux de surface"
$ws.Range("C202").Value = 66
$ws.Range("D202").Value = 1

# Reset style on the date cells so they match the default (unstyled) cells
# used throughout the rest of the sheet, instead of keeping the
# quote-prefix style variant that typing a leading apostrophe introduces.
$ws.Range("A199:A202").Style = "Normal"
